$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct tiny floating-point precision on the existing A6 timestamp
$ws.Range("A6").Value = 45865.20861030093

# Append the new row 7 with the latest sensor reading (06:00:22 run)
$ws.Range("A7").Value = 45865.25025472644
$ws.Range("A7").NumberFormat = $ws.Range("A6").NumberFormat
$ws.Range("B7").Value = 2025
$ws.Range("C7").Value = 30
$ws.Range("D7").Value = 13.11
$ws.Range("E7").Value = 90.69
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 2.42
$ws.Range("H7").Value = "ENE"
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = "06:00:22"
